$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.929.54"
$ws.Range("E2").Value = "  +0.27%  "

# Row 3
$ws.Range("D3").Value = "2.361.36"
$ws.Range("E3").Value = "  +2.02%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.86"
$ws.Range("E5").Value = "  -0.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.47"
$ws.Range("E6").Value = "  -0.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("E7").Value = "  -0.61%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("E9").Value = "  -1.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.88"
$ws.Range("E10").Value = "  -1.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("E11").Value = "  +0.11%  "

# Row 12
$ws.Range("E12").Value = "  +2.85%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.24"
$ws.Range("E13").Value = "  -3.91%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.71"
$ws.Range("E14").Value = "  -0.44%  "

# Row 15
$ws.Range("D15").Value = "2.730.78"
$ws.Range("E15").Value = "  +2.03%  "

# Row 16
$ws.Range("D16").Value = "2.363.33"
$ws.Range("E16").Value = "  +0.94%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.795"
$ws.Range("E17").Value = "  +0.83%  "

# Row 18
$ws.Range("D18").Value = "42.867.90"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.92"
$ws.Range("E19").Value = "  -2.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.25"
$ws.Range("E20").Value = "  +1.70%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0884"
$ws.Range("E21").Value = "  -0.82%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.82"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.93"
$ws.Range("E23").Value = "  -0.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("E24").Value = "  -2.91%  "

# Row 25
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.64"
$ws.Range("E27").Value = "  +0.98%  "

# Row 28
$ws.Range("E28").Value = "  +0.47%  "

# Row 29
$ws.Range("E29").Value = "  +1.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.51"
$ws.Range("E30").Value = "  -1.75%  "

# Row 31
$ws.Range("E31").Value = "  +0.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.03"
$ws.Range("E32").Value = "  +0.38%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.34"
$ws.Range("E33").Value = "  -2.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0724"
$ws.Range("E34").Value = "  +3.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.104"
$ws.Range("E35").Value = "  +3.99%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("E36").Value = "  +3.28%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.34"
$ws.Range("E37").Value = "  -2.79%  "

# Row 38
$ws.Range("E38").Value = "  -1.11%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.77"
$ws.Range("E39").Value = "  +1.28%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "120.83"
$ws.Range("E40").Value = "  -27.23%  "

# Row 41
$ws.Range("E41").Value = "  -0.76%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.54"
$ws.Range("E42").Value = "  +3.02%  "

# Row 43
$ws.Range("D43").Value = "1.932.51"
$ws.Range("E43").Value = "  +0.24%  "

# Row 44
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.14"
$ws.Range("E45").Value = "  +2.63%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.19"
$ws.Range("E46").Value = "  -9.43%  "

# Row 47
$ws.Range("E47").Value = "  -1.77%  "

# Row 48
$ws.Range("D48").Value = "2.589.23"
$ws.Range("E48").Value = "  +1.80%  "

# Row 49
$ws.Range("E49").Value = "  +2.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.72"
$ws.Range("E50").Value = "  -0.74%  "

# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.58"
$ws.Range("E51").Value = "  -3.33%  "
